# The "客単価" (customer unit price, column H) on the ABC分析_客構成 sheet
# had been accidentally left as a row total instead of a per-customer
# average. Recompute it by dividing the existing value by the customer
# count in column E (count_客構成) for each data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ABC分析_客構成")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row

for ($r = 3; $r -le $lastRow; $r++) {
    $countCell = $ws.Cells.Item($r, 5)
    $priceCell = $ws.Cells.Item($r, 8)

    $count = $countCell.Value()
    $price = $priceCell.Value()

    if ($count -and ($count -ne 0)) {
        $priceCell.Value = $price / $count
    }
}
